$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the _GoBack bookmark that sits at the end of the paragraph
#    "Three-Messengers-Verse: Removed Gloss styling on break character in
#    line 119" (last row of the corrections-log table, description cell).
#    Directly deleting a collapsed Range is unsafe in this runtime, so we
#    insert a clean replacement paragraph (formatting copied from the
#    still-bookmark-free paragraph that precedes it) immediately before the
#    original, move the text across, then delete the now-redundant original
#    paragraph (bookmark and all).
# ---------------------------------------------------------------------------
$targetText = "Three-Messengers-Verse: Removed Gloss styling on break character in line 119"

$rng = $d.Content
$rng.Find.Execute($targetText)
$origPara = $rng.Paragraphs(1)
$origStart = $origPara.Range.Start

# Locate the paragraph's 1-based index within $d.Paragraphs (".Previous"/
# ".Next" as bare properties are not reliable in this runtime, so resolve
# positionally instead).
$origIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Start -eq $origStart) {
        $origIndex = $i
        break
    }
}

$origPara.Range.InsertParagraphBefore()

# The new blank paragraph now occupies the original's old index; the
# original (still carrying its text + the _GoBack bookmark) has been pushed
# one paragraph later.
$newPara = $d.Paragraphs($origIndex)
$newPara.Range.InsertAfter($targetText)

$oldPara = $d.Paragraphs($origIndex + 1)
$oldPara.Range.Delete()

# ---------------------------------------------------------------------------
# 2) Append two new rows to the (only) table, recording the latest
#    corrections made by Steffi Delcourt on 10.30.24.
# ---------------------------------------------------------------------------
$t = $d.Tables(1)

$row1 = $t.Rows.Add()
$i1 = $row1.Index
$t.Cell($i1, 1).Range.Text = "10.30.24"
$t.Cell($i1, 2).Range.Text = "12-Deaths-Warning-Verse: added a missing sidegloss " + [char]8220 + "again" + [char]8221 + " to line 52"
$t.Cell($i1, 3).Range.Text = "Steffi Delcourt"

$row2 = $t.Rows.Add()
$i2 = $row2.Index
$t.Cell($i2, 1).Range.Text = "10.30.24"
$t.Cell($i2, 2).Range.Text = "27-Ressoning-Verse: found and deleted forced line breaks that were incorrectly styled with " + [char]8220 + "line number" + [char]8221 + " CS " + [char]8211 + " 4 instances corrected"
$t.Cell($i2, 3).Range.Text = "Steffi Delcourt"
